$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: CSC103 -> MAT141
$ws.Range("B2").Value = "MAT141"
$ws.Range("C2").Value = "store/materials/MAT141/conditional probability.pdf"

# Update row 3: MAT141 -> MAT111
$ws.Range("B3").Value = "MAT111"
$ws.Range("C3").Value = "store/materials/MAT111/10 Usability Heuristics for User Interface Design_1622399977365.pdf"
